$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data rows 2-16: column A (word), column B (predicted word), column C (score)
$data = @(
    @("<they>",   "<they>", 40),
    @("<her>",    "<when>", 39),
    @("<his>",    "<is>",   43),
    @("<had>",    "<can>",  32),
    @("<find>",   "<on>",   48),
    @("<of>",     "<of>",   38),
    @("<juliet>", "<oil>",  37),
    @("<word>",   "<when>", 34),
    @("<paste>",  "<be>",   39),
    @("<make>",   "<make>", 41),
    @("<yes>",    "<you>",  38),
    @("<part>",   "<oil>",  34),
    @("<many>",   "<many>", 33),
    @("<lima>",   "<more>", 31),
    @("<their>",  "<their>", 22)
)

$row = 2
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $row++
}

$wb.Save()
